# Update '想去人数' (Column F) view/interest counts across all sheets
# per the gh-pages data regeneration (commit 456a3b4).

$wb = $excel.ActiveWorkbook

# Sheet: 展览
$ws = $wb.Worksheets.Item("展览")
$ws.Cells.Item(3, 6).Value = 637  # F3: 635 -> 637
$ws.Cells.Item(4, 6).Value = 2820  # F4: 2813 -> 2820
$ws.Cells.Item(8, 6).Value = 278  # F8: 276 -> 278
$ws.Cells.Item(9, 6).Value = 6534  # F9: 6522 -> 6534
$ws.Cells.Item(13, 6).Value = 5089  # F13: 5083 -> 5089
$ws.Cells.Item(15, 6).Value = 565  # F15: 563 -> 565
$ws.Cells.Item(16, 6).Value = 2690  # F16: 2679 -> 2690
$ws.Cells.Item(17, 6).Value = 1369  # F17: 1368 -> 1369
$ws.Cells.Item(18, 6).Value = 1533  # F18: 1531 -> 1533
$ws.Cells.Item(19, 6).Value = 1242  # F19: 1241 -> 1242
$ws.Cells.Item(20, 6).Value = 329  # F20: 328 -> 329
$ws.Cells.Item(21, 6).Value = 131  # F21: 130 -> 131
$ws.Cells.Item(22, 6).Value = 146  # F22: 143 -> 146
$ws.Cells.Item(23, 6).Value = 1114  # F23: 1110 -> 1114
$ws.Cells.Item(24, 6).Value = 262  # F24: 259 -> 262
$ws.Cells.Item(25, 6).Value = 556  # F25: 552 -> 556
$ws.Cells.Item(26, 6).Value = 1401  # F26: 1400 -> 1401
$ws.Cells.Item(28, 6).Value = 2124  # F28: 2122 -> 2124
$ws.Cells.Item(29, 6).Value = 604  # F29: 599 -> 604
$ws.Cells.Item(30, 6).Value = 51  # F30: 49 -> 51
$ws.Cells.Item(31, 6).Value = 48  # F31: 44 -> 48
$ws.Cells.Item(32, 6).Value = 119  # F32: 117 -> 119
$ws.Cells.Item(33, 6).Value = 271  # F33: 267 -> 271
$ws.Cells.Item(34, 6).Value = 1545  # F34: 1543 -> 1545
$ws.Cells.Item(35, 6).Value = 9  # F35: 7 -> 9
$ws.Cells.Item(36, 6).Value = 8  # F36: 7 -> 8
$ws.Cells.Item(38, 6).Value = 1103  # F38: 1101 -> 1103
$ws.Cells.Item(41, 6).Value = 2332  # F41: 2329 -> 2332
$ws.Cells.Item(42, 6).Value = 2595  # F42: 2593 -> 2595
$ws.Cells.Item(44, 6).Value = 155  # F44: 152 -> 155
$ws.Cells.Item(46, 6).Value = 284  # F46: 282 -> 284
$ws.Cells.Item(48, 6).Value = 117  # F48: 115 -> 117
$ws.Cells.Item(49, 6).Value = 391  # F49: 390 -> 391

# Sheet: 演出
$ws = $wb.Worksheets.Item("演出")
$ws.Cells.Item(9, 6).Value = 333  # F9: 332 -> 333
$ws.Cells.Item(11, 6).Value = 170  # F11: 169 -> 170
$ws.Cells.Item(12, 6).Value = 100  # F12: 99 -> 100
$ws.Cells.Item(13, 6).Value = 205  # F13: 204 -> 205
$ws.Cells.Item(17, 6).Value = 167  # F17: 166 -> 167
$ws.Cells.Item(18, 6).Value = 46  # F18: 44 -> 46
$ws.Cells.Item(27, 6).Value = 428  # F27: 427 -> 428

# Sheet: 本地生活
$ws = $wb.Worksheets.Item("本地生活")
$ws.Cells.Item(4, 6).Value = 522  # F4: 520 -> 522
$ws.Cells.Item(6, 6).Value = 1710  # F6: 1711 -> 1710
$ws.Cells.Item(8, 6).Value = 1553  # F8: 1550 -> 1553
$ws.Cells.Item(9, 6).Value = 1827  # F9: 1826 -> 1827
$ws.Cells.Item(10, 6).Value = 2594  # F10: 2593 -> 2594
$ws.Cells.Item(11, 6).Value = 903  # F11: 899 -> 903
$ws.Cells.Item(12, 6).Value = 787  # F12: 785 -> 787
$ws.Cells.Item(14, 6).Value = 161  # F14: 159 -> 161

# Sheet: 全部类型
$ws = $wb.Worksheets.Item("全部类型")
$ws.Cells.Item(3, 6).Value = 522  # F3: 520 -> 522
$ws.Cells.Item(4, 6).Value = 1710  # F4: 1711 -> 1710
$ws.Cells.Item(5, 6).Value = 637  # F5: 635 -> 637
$ws.Cells.Item(6, 6).Value = 2820  # F6: 2813 -> 2820
$ws.Cells.Item(8, 6).Value = 1553  # F8: 1550 -> 1553
$ws.Cells.Item(9, 6).Value = 278  # F9: 276 -> 278
$ws.Cells.Item(10, 6).Value = 6534  # F10: 6522 -> 6534
$ws.Cells.Item(11, 6).Value = 903  # F11: 899 -> 903
$ws.Cells.Item(12, 6).Value = 787  # F12: 785 -> 787
$ws.Cells.Item(13, 6).Value = 5089  # F13: 5083 -> 5089
$ws.Cells.Item(15, 6).Value = 565  # F15: 563 -> 565
$ws.Cells.Item(16, 6).Value = 2690  # F16: 2679 -> 2690
$ws.Cells.Item(17, 6).Value = 1369  # F17: 1368 -> 1369
$ws.Cells.Item(18, 6).Value = 1242  # F18: 1241 -> 1242
$ws.Cells.Item(19, 6).Value = 329  # F19: 328 -> 329
$ws.Cells.Item(21, 6).Value = 131  # F21: 130 -> 131
$ws.Cells.Item(22, 6).Value = 146  # F22: 143 -> 146
$ws.Cells.Item(23, 6).Value = 333  # F23: 332 -> 333
$ws.Cells.Item(24, 6).Value = 1114  # F24: 1110 -> 1114
$ws.Cells.Item(25, 6).Value = 262  # F25: 259 -> 262
$ws.Cells.Item(26, 6).Value = 100  # F26: 99 -> 100
$ws.Cells.Item(27, 6).Value = 161  # F27: 159 -> 161
$ws.Cells.Item(28, 6).Value = 556  # F28: 553 -> 556
$ws.Cells.Item(29, 6).Value = 1401  # F29: 1400 -> 1401
$ws.Cells.Item(31, 6).Value = 2124  # F31: 2122 -> 2124
$ws.Cells.Item(32, 6).Value = 604  # F32: 599 -> 604
$ws.Cells.Item(33, 6).Value = 51  # F33: 49 -> 51
$ws.Cells.Item(34, 6).Value = 167  # F34: 166 -> 167
$ws.Cells.Item(35, 6).Value = 48  # F35: 44 -> 48
$ws.Cells.Item(36, 6).Value = 271  # F36: 267 -> 271
$ws.Cells.Item(37, 6).Value = 46  # F37: 44 -> 46
$ws.Cells.Item(38, 6).Value = 1545  # F38: 1543 -> 1545
$ws.Cells.Item(39, 6).Value = 9  # F39: 7 -> 9
$ws.Cells.Item(40, 6).Value = 1103  # F40: 1101 -> 1103
$ws.Cells.Item(44, 6).Value = 2332  # F44: 2329 -> 2332
$ws.Cells.Item(45, 6).Value = 2595  # F45: 2593 -> 2595
$ws.Cells.Item(46, 6).Value = 155  # F46: 152 -> 155
$ws.Cells.Item(49, 6).Value = 391  # F49: 390 -> 391

